$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '58.086.34'
$ws.Range('E2').Value = '  -0.77%  '
$ws.Range('D3').Value = '2.284.17'
$ws.Range('E3').Value = '  +0.36%  '
$ws.Range('E4').Value = '  +0.05%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '534.72'
$ws.Range('D5').ClearFormats()
$ws.Range('E5').Value = '  -1.98%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '131.26'
$ws.Range('D6').ClearFormats()
$ws.Range('E6').Value = '  +0.60%  '
$ws.Range('E7').Value = '  +0.05%  '
$ws.Range('E8').Value = '  +3.54%  '
$ws.Range('D9').Value = '2.284.39'
$ws.Range('E9').Value = '  +0.48%  '
$ws.Range('E10').Value = '  -1.47%  '
$ws.Range('E11').Value = '  -0.44%  '
$ws.Range('E12').Value = '  +0.66%  '
$ws.Range('E13').Value = '  -0.51%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '23.50'
$ws.Range('D14').ClearFormats()
$ws.Range('E14').Value = '  -0.65%  '
$ws.Range('D15').Value = '2.691.43'
$ws.Range('E15').Value = '  +0.37%  '
$ws.Range('D16').Value = '58.037.16'
$ws.Range('E16').Value = '  -0.79%  '
$ws.Range('E17').Value = '  -0.69%  '
$ws.Range('D18').Value = '2.295.68'
$ws.Range('E18').Value = '  +0.93%  '
$ws.Range('E19').Value = '  -1.38%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '4.18'
$ws.Range('D20').ClearFormats()
$ws.Range('E20').Value = '  -2.50%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '313.18'
$ws.Range('D21').ClearFormats()
$ws.Range('E21').Value = '  -0.30%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '6.44'
$ws.Range('D22').ClearFormats()
$ws.Range('E22').Value = '  +0.18%  '
$ws.Range('E23').Value = '  -0.01%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '63.08'
$ws.Range('D24').ClearFormats()
$ws.Range('E24').Value = '  +0.35%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '0.168'
$ws.Range('D25').ClearFormats()
$ws.Range('E25').Value = '  -1.14%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '0.998'
$ws.Range('D26').ClearFormats()
$ws.Range('E26').Value = '  -0.20%  '
$ws.Range('E27').Value = '  -1.76%  '
$ws.Range('E28').Value = '  -2.99%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '170.91'
$ws.Range('D29').ClearFormats()
$ws.Range('E29').Value = '  +0.11%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '1.71'
$ws.Range('D30').ClearFormats()
$ws.Range('E30').Value = '  -1.97%  '
$ws.Range('D31').Value = '0.0₃0721'
$ws.Range('E31').Value = '  +0.22%  '
$ws.Range('B32').Value = 'Aptos'
$ws.Range('C32').Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '5.76'
$ws.Range('D32').ClearFormats()
$ws.Range('E32').Value = '  +0.05%  '
$ws.Range('B33').Value = 'SuiNetwork'
$ws.Range('C33').Value = 'https://coinranking.com/coin/3xJluUMvp+suinetwork-sui'
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '1.07'
$ws.Range('D33').ClearFormats()
$ws.Range('E33').Value = '  -0.45%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '0.381'
$ws.Range('D34').ClearFormats()
$ws.Range('E34').Value = '  -0.36%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '17.82'
$ws.Range('D36').ClearFormats()
$ws.Range('E36').Value = '  +0.48%  '
$ws.Range('E37').Value = '  -0.05%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '1.23'
$ws.Range('D38').ClearFormats()
$ws.Range('E38').Value = '  -1.27%  '
$ws.Range('E39').Value = '  -1.12%  '
$ws.Range('E40').Value = '  -1.49%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '288.51'
$ws.Range('D41').ClearFormats()
$ws.Range('E41').Value = '  -4.04%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '139.78'
$ws.Range('D42').ClearFormats()
$ws.Range('E42').Value = '  -0.53%  '
$ws.Range('E43').Value = '  -0.34%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '0.0951'
$ws.Range('D44').ClearFormats()
$ws.Range('E44').Value = '  +0.58%  '
$ws.Range('E45').Value = '  -0.59%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '0.552'
$ws.Range('D46').ClearFormats()
$ws.Range('E46').Value = '  +0.53%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '18.06'
$ws.Range('D47').ClearFormats()
$ws.Range('E47').Value = '  -1.38%  '
$ws.Range('E48').Value = '  -1.53%  '
$ws.Range('E49').Value = '  -0.53%  '
$ws.Range('E50').Value = '  +0.02%  '
$ws.Range('E51').Value = '  +1.06%  '
